$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new purchase row for 2025-05-01 as plain text/values,
# matching the pattern used for the other manually-entered rows
# (A column stored as text date strings, B/C/D as numbers).
$ws.Range("A18").NumberFormat = "@"
$ws.Range("A18").Value = "05/01/2025"
$ws.Range("A18").Style = "Normal"
$ws.Range("B18").Value = 524.0279999999984
$ws.Range("C18").Value = 0.09541474883021546
$ws.Range("D18").Value = 50
